$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name / text updates ---
$ws.Range("A5").Value = "Major Gen. Brighton"
$ws.Range("A6").Value = "Governor Cartwright"
$ws.Range("A11").Value = "Vice President Ramsey"
$ws.Range("A13").Value = "Dr. Jacob Neumann"
$ws.Range("A14").Value = "Layabout"

# --- Row 15 (Researcher) ---
$ws.Range("G15").Value = 0
$ws.Range("J15").Value = "+"

# --- Row 16 (Tax Collector) ---
$ws.Range("B16").Value = "1N1R"
$ws.Range("K16").Value = 2

# --- Row 17 (Ambassador) ---
$ws.Range("B17").Value = "1X"

# --- Row 18 (Philanthropist -> Wealthy Donor) ---
$ws.Range("A18").Value = "Wealthy Donor"
$ws.Range("K18").Value = 2
$ws.Range("M18").Value = "+"

# --- Row 19 (Merchant -> Captain of Industry) ---
$ws.Range("A19").Value = "Captain of Industry"
$ws.Range("B19").Value = "1N3X"

# --- Selection change (active cell) ---
$ws.Range("B18").Select()
